$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the student-vaccination header row (USN, Name, Age, Phone, Vaccine_Dose)
$ws.Range("A1").Value = "USN"
$ws.Range("B1").Value = "Name"
$ws.Range("C1").Value = "Age"
$ws.Range("D1").Value = "Phone"
$ws.Range("E1").Value = "Vaccine_Dose"

# Widen column E to fit its header text, just like Excel does when a user
# double-clicks the column border after typing a long heading
$ws.Columns("E:E").AutoFit()

# Leave the selection on the first empty cell after the header, matching
# where the cursor lands after tabbing across the row
$ws.Range("F1").Select()
